$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 8086
$ws.Range("F5").Value = 30658
$ws.Range("F7").Value = 605
$ws.Range("F8").Value = 698
$ws.Range("F9").Value = 459
$ws.Range("F11").Value = 445
$ws.Range("F12").Value = 798
$ws.Range("F13").Value = 52
$ws.Range("F14").Value = 604
$ws.Range("F15").Value = 386
$ws.Range("F17").Value = 556
$ws.Range("F19").Value = 412
$ws.Range("F22").Value = 85
$ws.Range("F23").Value = 701
$ws.Range("F24").Value = 2333
$ws.Range("F25").Value = 826
$ws.Range("F26").Value = 66
$ws.Range("F29").Value = 633
$ws.Range("F30").Value = 1073

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 342

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 521

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 521
$ws.Range("F3").Value = 8086
$ws.Range("F7").Value = 30658
$ws.Range("F9").Value = 605
$ws.Range("F10").Value = 698
$ws.Range("F11").Value = 459
$ws.Range("F14").Value = 445
$ws.Range("F15").Value = 342
$ws.Range("F18").Value = 798
$ws.Range("F19").Value = 52
$ws.Range("F20").Value = 604
$ws.Range("F21").Value = 386
$ws.Range("F27").Value = 556
$ws.Range("F29").Value = 412
$ws.Range("F32").Value = 85
$ws.Range("F33").Value = 701
$ws.Range("F34").Value = 2333
$ws.Range("F35").Value = 826
$ws.Range("F36").Value = 66
$ws.Range("F40").Value = 633
$ws.Range("F41").Value = 1073
